# Add "Stricter Marketing policy" variable (Model 4) to the workbook.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data for the model")
$wsExpl = $wb.Worksheets.Item("explanation of variables")

# --- 1. Add new header "Stricter marketing policy" in column I of the data sheet ---
$wsData.Range("I1").Value = "Stricter marketing policy"

# I column values for rows 2..81 (0/1 indicator), 1 at rows 9,19,34,49,59
$iValues = @(0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 9).Value = $iValues[$i]
}

# widen column I to fit
$wsData.Columns.Item(9).ColumnWidth = 21.08984375

# --- 2. Add explanatory row 10 on "explanation of variables" sheet ---
$wsExpl.Range("A10").Value = "Stricter marketing policy"
$wsExpl.Range("B10").Value = "Restrictions on alcohol marketing"
$wsExpl.Range("A10").Style = $wsExpl.Range("A9").Style
$wsExpl.Range("B10").Style = $wsExpl.Range("B9").Style

# The table "Tabela1" on this sheet should auto-expand to include the new row.
$tbl = $wsExpl.ListObjects.Item("Tabela1")
$tbl.Resize($wsExpl.Range("A1:B10"))

# Column A width adjustment on explanation sheet
$wsExpl.Columns.Item(1).ColumnWidth = 20.26953125
